$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the two algorithm/label entries in column F (rows 3 and 4)
$ws.Range("F3").Value = "TA Ciclica Aleatoria Ext"
$ws.Range("F4").Value = " TA Decreciente (propuesta)"

# Match the author's last selected cell when the workbook was saved
[void]$ws.Range("F7").Select()
